# Apply the "New DB connection working." update:
#  - Row 5: Lastname "DEAN" -> "TESET", Firstname "TARA" -> "TEST"
#  - Active selection moves from G8 to C5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "TEST"
$ws.Range("C5").Value = "TESET"

$ws.Range("C5").Select()
